$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.358.52"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.712.15"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.56"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5297"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06703"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2665"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.90"
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07676"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.512"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "1.947.11"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "1.712.72"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5833"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "0.0₅8232"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.15"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "27.344.81"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.03"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.012"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.90"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.690"
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.240"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.33"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05365"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.483"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.435"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.637"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.872"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9504"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.398"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5854"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01635"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "1.094.11"
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.792"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8413"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.97"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.853.98"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.77"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4538"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.081"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  -0.29%  "
